# Auto-generated Excel COM-interop script applying the scheduled-runner Sheets update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 129 (ALC)
$ws.Range("H129").Value = 1023.5455
$ws.Range("I129").Value = 481.75
$ws.Range("J129").Value = 1143.9445
$ws.Range("K129").Value = 1445.25
$ws.Range("L129").Value = 3431.8335
$ws.Range("M129").Value = 3554.75
$ws.Range("N129").Value = -13431.8335

# Row 133 (ALC)
$ws.Range("H133").Value = 59600
$ws.Range("J133").Value = 59600
$ws.Range("L133").Value = 59600
$ws.Range("N133").Value = -69720

# Row 137 (ALC)
$ws.Range("H137").Value = 2107.25
$ws.Range("I137").Value = 2633.3333
$ws.Range("J137").Value = 1791.6
$ws.Range("K137").Value = 7899.999899999999
$ws.Range("L137").Value = 5374.799999999999
$ws.Range("M137").Value = -5349.999899999999
$ws.Range("N137").Value = -10474.8

# Row 138 (ALC)
$ws.Range("H138").Value = 3651.2458
$ws.Range("I138").Value = 1566.9062
$ws.Range("J138").Value = 5951.207
$ws.Range("K138").Value = 4700.7186
$ws.Range("L138").Value = 17853.621
$ws.Range("M138").Value = 439.2813999999998
$ws.Range("N138").Value = -28133.621

$ws = $wb.Worksheets.Item("ARM")
# Row 37 (ARM)
$ws.Range("H37").Value = 11038
$ws.Range("J37").Value = 11038
$ws.Range("L37").Value = 11038
$ws.Range("N37").Value = -11584

# Row 55 (ARM)
$ws.Range("H55").Value = 34399.125
$ws.Range("J55").Value = 34399.125
$ws.Range("L55").Value = 34399.125
$ws.Range("N55").Value = -35029.125

# Row 80 (ARM)
$ws.Range("H80").Value = 20375.25
$ws.Range("J80").Value = 20375.25
$ws.Range("L80").Value = 20375.25
$ws.Range("N80").Value = -22371.25

# Row 83 (ARM)
$ws.Range("H83").Value = 20375.25
$ws.Range("J83").Value = 20375.25
$ws.Range("L83").Value = 61125.75
$ws.Range("N83").Value = -71109.75

# Row 132 (ARM)
$ws.Range("H132").Value = 21299670
$ws.Range("I132").Value = 29413146
$ws.Range("K132").Value = 88239438
$ws.Range("M132").Value = -88236908

# Row 133 (ARM)
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -35060

$ws = $wb.Worksheets.Item("BSM")
# Row 81 (BSM)
$ws.Range("H81").Value = 25822.223
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 25822.223
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 25822.223
$ws.Range("N81").Value = -27944.223
$ws.Range("M81").ClearContents()

# Row 84 (BSM)
$ws.Range("H84").Value = 25822.223
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 25822.223
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 77466.66900000001
$ws.Range("N84").Value = -88074.66900000001
$ws.Range("M84").ClearContents()

# Row 132 (BSM)
$ws.Range("H132").Value = 54390
$ws.Range("J132").Value = 54390
$ws.Range("L132").Value = 54390
$ws.Range("N132").Value = -64510

# Row 134 (BSM)
$ws.Range("H134").Value = 12240.211
$ws.Range("I134").Value = 4985.625
$ws.Range("K134").Value = 14956.875
$ws.Range("M134").Value = -12421.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 2787522.5
$ws.Range("I31").Value = 3248444
$ws.Range("K31").Value = 3248444
$ws.Range("M31").Value = -3248149

# Row 34 (CRP)
$ws.Range("H34").Value = 2787522.5
$ws.Range("I34").Value = 3248444
$ws.Range("K34").Value = 3248444
$ws.Range("M34").Value = -3248242

# Row 86 (CRP)
$ws.Range("H86").Value = 25067.334
$ws.Range("I86").Value = 26250
$ws.Range("K86").Value = 26250
$ws.Range("M86").Value = -25127

# Row 89 (CRP)
$ws.Range("H89").Value = 25067.334
$ws.Range("I89").Value = 26250
$ws.Range("K89").Value = 131250
$ws.Range("M89").Value = -125634

# Row 132 (CRP)
$ws.Range("H132").Value = 4653236
$ws.Range("I132").Value = 1906.9
$ws.Range("J132").Value = 15387072
$ws.Range("K132").Value = 5720.700000000001
$ws.Range("L132").Value = 46161216
$ws.Range("M132").Value = -3190.700000000001
$ws.Range("N132").Value = -46166276

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 1156.8334
$ws.Range("I5").Value = 875.7778
$ws.Range("K5").Value = 2627.3334
$ws.Range("M5").Value = -2515.3334

# Row 122 (CUL)
$ws.Range("H122").Value = 830.8
$ws.Range("I122").Value = 325
$ws.Range("K122").Value = 2925
$ws.Range("M122").Value = -475

# Row 129 (CUL)
$ws.Range("H129").Value = 13334718
$ws.Range("I129").Value = 955
$ws.Range("K129").Value = 2865
$ws.Range("M129").Value = 2135

# Row 131 (CUL)
$ws.Range("H131").Value = 5748030
$ws.Range("J131").Value = 6945397.5
$ws.Range("L131").Value = 20836192.5
$ws.Range("N131").Value = -20846272.5

# Row 132 (CUL)
$ws.Range("H132").Value = 1718.9445
$ws.Range("J132").Value = 2338.1
$ws.Range("L132").Value = 21042.9
$ws.Range("N132").Value = -26102.9

# Row 135 (CUL)
$ws.Range("H135").Value = 1156.8334
$ws.Range("I135").Value = 875.7778
$ws.Range("K135").Value = 7882.000199999999
$ws.Range("M135").Value = -5347.000199999999

# Row 139 (CUL)
$ws.Range("H139").Value = 2320
$ws.Range("I139").Value = 1982.8572
$ws.Range("K139").Value = 5948.571599999999
$ws.Range("M139").Value = -808.5715999999993

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 1793207.9
$ws.Range("I70").Value = 2724038.2
$ws.Range("J70").Value = 9116.666999999999
$ws.Range("K70").Value = 2724038.2
$ws.Range("L70").Value = 9116.666999999999
$ws.Range("M70").Value = -2723768.2
$ws.Range("N70").Value = -9656.666999999999

# Row 73 (GSM)
$ws.Range("H73").Value = 1793207.9
$ws.Range("I73").Value = 2724038.2
$ws.Range("J73").Value = 9116.666999999999
$ws.Range("K73").Value = 2724038.2
$ws.Range("L73").Value = 9116.666999999999
$ws.Range("M73").Value = -2723102.2
$ws.Range("N73").Value = -10988.667

# Row 80 (GSM)
$ws.Range("H80").Value = 3843.3333
$ws.Range("I80").Value = 2205
$ws.Range("J80").Value = 4171
$ws.Range("K80").Value = 2205
$ws.Range("L80").Value = 4171
$ws.Range("M80").Value = -1207
$ws.Range("N80").Value = -6167

# Row 83 (GSM)
$ws.Range("H83").Value = 3843.3333
$ws.Range("I83").Value = 2205
$ws.Range("J83").Value = 4171
$ws.Range("K83").Value = 11025
$ws.Range("L83").Value = 20855
$ws.Range("M83").Value = -6033
$ws.Range("N83").Value = -30839

# Row 126 (GSM)
$ws.Range("H126").Value = 9811534
$ws.Range("I126").Value = 9637.25
$ws.Range("J126").Value = 33336086
$ws.Range("K126").Value = 28911.75
$ws.Range("L126").Value = 100008258
$ws.Range("M126").Value = -26441.75
$ws.Range("N126").Value = -100013198

# Row 132 (GSM)
$ws.Range("H132").Value = 128048.5
$ws.Range("I132").Value = 2515.5
$ws.Range("J132").Value = 253581.5
$ws.Range("K132").Value = 7546.5
$ws.Range("L132").Value = 760744.5
$ws.Range("M132").Value = -5016.5
$ws.Range("N132").Value = -765804.5

# Row 133 (GSM)
$ws.Range("H133").Value = 48000
$ws.Range("J133").Value = 48000
$ws.Range("L133").Value = 48000
$ws.Range("N133").Value = -58120

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (LTW)
$ws.Range("H82").Value = 2147.1333
$ws.Range("I82").Value = 1978
$ws.Range("J82").Value = 2400.8333
$ws.Range("K82").Value = 1978
$ws.Range("L82").Value = 2400.8333
$ws.Range("M82").Value = -1617
$ws.Range("N82").Value = -3122.8333

# Row 85 (LTW)
$ws.Range("H85").Value = 2147.1333
$ws.Range("I85").Value = 1978
$ws.Range("J85").Value = 2400.8333
$ws.Range("K85").Value = 1978
$ws.Range("L85").Value = 2400.8333
$ws.Range("M85").Value = -730
$ws.Range("N85").Value = -4896.8333

# Row 100 (LTW)
$ws.Range("H100").Value = 2435.9246
$ws.Range("I100").Value = 1051.5
$ws.Range("J100").Value = 2548.9387
$ws.Range("K100").Value = 1051.5
$ws.Range("L100").Value = 2548.9387
$ws.Range("M100").Value = -510.5
$ws.Range("N100").Value = -3630.9387

# Row 132 (LTW)
$ws.Range("H132").Value = 10420252
$ws.Range("I132").Value = 41668120
$ws.Range("J132").Value = 4296.3335
$ws.Range("K132").Value = 125004360
$ws.Range("L132").Value = 12889.0005
$ws.Range("M132").Value = -125001830
$ws.Range("N132").Value = -17949.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 210813170
$ws.Range("I132").Value = 321429340
$ws.Range("J132").Value = 17234886
$ws.Range("K132").Value = 964288020
$ws.Range("L132").Value = 51704658
$ws.Range("M132").Value = -964285490
$ws.Range("N132").Value = -51709718

# Row 133 (WVR)
$ws.Range("H133").Value = 46000
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
